# Insert a new weekly price-record row at row 12 (right after the first
# 10 data rows), pushing all subsequent rows down by one. Fill the new
# row with the new week's data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12..97 down to 13..98, creating a blank row 12.
$ws.Rows("12:12").Insert()

# Populate the new row 12 with the new data point.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44532
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 100112031
$ws.Range("G12").Value = "Poroto verde"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 125
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 25000
$ws.Range("N12").Value = "$/malla 25 kilos"
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1000
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
